$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
    $ws.Range($range).Style = "Normal"
}

# Row 2 (Bitcoin)
Set-TextValue "D2" "62.595.33"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3 (Ethereum)
Set-TextValue "D3" "2.451.21"
$ws.Range("E3").Value = "  -0.86%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 (BNB)
Set-TextValue "D5" "570.14"
$ws.Range("E5").Value = "  -1.24%  "

# Row 6 (Solana)
Set-TextValue "D6" "145.75"
$ws.Range("E6").Value = "  -0.64%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 (XRP)
$ws.Range("E8").Value = "  -2.07%  "

# Row 9 (Dogecoin)
$ws.Range("E9").Value = "  -1.15%  "

# Row 10 (TRON)
$ws.Range("E10").Value = "  -0.29%  "

# Row 11 (Toncoin)
$ws.Range("E11").Value = "  -2.55%  "

# Row 12 (Cardano)
$ws.Range("E12").Value = "  -1.64%  "

# Row 13 (Avalanche)
Set-TextValue "D13" "28.49"
$ws.Range("E13").Value = "  -1.72%  "

# Row 14 (ShibaInu)
$ws.Range("E14").Value = "  -3.20%  "

# Row 15 (WrappedliquidstakedEther2.0)
Set-TextValue "D15" "2.895.50"

# Row 16 (WrappedBTC)
Set-TextValue "D16" "62.402.00"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17 (WrappedEther)
Set-TextValue "D17" "2.452.52"
$ws.Range("E17").Value = "  -0.90%  "

# Row 18 (Uniswap)
Set-TextValue "D18" "7.65"
$ws.Range("E18").Value = "  -6.51%  "

# Row 19 (Chainlink)
Set-TextValue "D19" "10.67"
$ws.Range("E19").Value = "  -3.17%  "

# Row 20 (BitcoinCash)
Set-TextValue "D20" "320.53"
$ws.Range("E20").Value = "  -2.81%  "

# Row 21 (Polkadot)
Set-TextValue "D21" "4.12"
$ws.Range("E21").Value = "  -0.43%  "

# Row 22 (SuiNetwork)
$ws.Range("E22").Value = "  -1.11%  "

# Row 23 (Dai)
Set-TextValue "D23" "1.00"
$ws.Range("E23").Value = "  +0.05%  "

# Row 24 (Aptos)
Set-TextValue "D24" "9.87"
$ws.Range("E24").Value = "  +2.45%  "

# Row 25 (Litecoin)
Set-TextValue "D25" "64.64"
$ws.Range("E25").Value = "  -2.54%  "

# Row 26 (Bittensor)
Set-TextValue "D26" "644.11"
$ws.Range("E26").Value = "  -3.30%  "

# Row 27 (WrappedeETH)
Set-TextValue "D27" "2.571.57"
$ws.Range("E27").Value = "  -0.76%  "

# Row 28 (Binance-PegBSC-USD)
Set-TextValue "D28" "0.998"
$ws.Range("E28").Value = "  -0.35%  "

# Row 29 (PEPE)
Set-TextValue "D29" "0.0₃0943"
$ws.Range("E29").Value = "  -4.43%  "

# Row 30 (Fetch.AI)
$ws.Range("E30").Value = "  -3.58%  "

# Row 31 (InternetComputer(DFINITY))
Set-TextValue "D31" "7.77"
$ws.Range("E31").Value = "  -3.83%  "

# Row 32 (PancakeSwap)
$ws.Range("E32").Value = "  -3.50%  "

# Row 33 (Kaspa)
Set-TextValue "D33" "0.132"
$ws.Range("E33").Value = "  -0.79%  "

# Row 34 (FirstDigitalUSD)
$ws.Range("E34").Value = "  -0.04%  "

# Row 35 (ImmutableX)
$ws.Range("E35").Value = "  -4.01%  "

# Row 36 (Monero)
Set-TextValue "D36" "151.75"
$ws.Range("E36").Value = "  -0.61%  "

# Row 37 (NEARProtocol)
$ws.Range("E37").Value = "  -3.91%  "

# Row 38 (EthereumClassic)
Set-TextValue "D38" "18.47"
$ws.Range("E38").Value = "  -1.64%  "

# Row 39 (PolygonEcosystemToken)
$ws.Range("E39").Value = "  -2.56%  "

# Row 40 (RenderToken)
Set-TextValue "D40" "5.27"
$ws.Range("E40").Value = "  -2.98%  "

# Row 41 (dogwifhat)
$ws.Range("E41").Value = "  -4.19%  "

# Row 42 (Stacks)
$ws.Range("E42").Value = "  -4.24%  "

# Row 43 (USDe)
$ws.Range("E43").Value = "  +0.04%  "

# Row 44 (BabyDogeCoin)
Set-TextValue "D44" "0.0₆0305"
$ws.Range("E44").Value = "  +0.30%  "

# Row 45 (Aave)
Set-TextValue "D45" "151.99"
$ws.Range("E45").Value = "  +0.18%  "

# Row 46 (WhiteBITCoin)
Set-TextValue "D46" "15.39"
$ws.Range("E46").Value = "  +1.62%  "

# Row 47 (Filecoin)
$ws.Range("E47").Value = "  -2.48%  "

# Row 48 (Mantle)
$ws.Range("E48").Value = "  -0.92%  "

# Row 49 (InjectiveProtocol)
$ws.Range("E49").Value = "  -4.51%  "

# Row 50 (Hedera)
$ws.Range("E50").Value = "  -2.54%  "

# Row 51 (Stellar)
$ws.Range("E51").Value = "  -2.20%  "
